# Actualización automática de tasas-transfi.xlsx
# Updates the "Conversión del día" note on Hoja1 and the N/O rate table on tasas.

$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: refresh Binance conversion figures ---------------------------
$hoja1 = $wb.Worksheets.Item("Hoja1")

$nuevoTexto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.56 = 25954.1 pesos`n✅ 25954.1 pesos = 6.56 = 977.79 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$hoja1.Range("A1").Value = $nuevoTexto

# --- tasas!N10/O10/N12/O12: refresh the Binance/transfi rate table ----------
$tasas = $wb.Worksheets.Item("tasas")

$tasas.Range("N10").Value = 152.5
$tasas.Range("O10").Value = 3958
$tasas.Range("N12").Value = 3955
$tasas.Range("O12").Value = 149
